$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 headers ---
$ws.Range("C1").Value = "Object 1"
$ws.Range("D1").Value = "Object 2"
$ws.Range("E1").Value = "Object 3"
$ws.Range("F1").Value = "Object 4"
$ws.Range("G1").Value = "Object 5"
$ws.Range("H1").Value = "Object 6"
$ws.Range("I1").Value = "Area in sq m"
$ws.Range("J1").Value = "Area in sq cm"

# --- Row 2: move the "model" area figures from C2/D2 into I2/J2 ---
$ws.Range("C2:D2").ClearContents()
$ws.Range("I2").Value = 5.4051768899999999
$ws.Range("J2").Formula = "=I2*100*100"

# --- Row 3: coral measurements across 6 objects ---
$ws.Range("B3").Value = 6
$ws.Range("C3").Value = 0.16218480199999999
$ws.Range("D3").Value = 0.0170914707
$ws.Range("E3").Value = 0.00843721478
$ws.Range("F3").Value = 0.00685544447
$ws.Range("G3").Value = 0.00551765668
$ws.Range("H3").Value = 0.00354007093
$ws.Range("I3").Formula = "=SUM(C3:H3)"

# --- Row 4: disease measurements across a subset of objects ---
$ws.Range("C4").Value = 0.00296883294
$ws.Range("D4").Value = 0.00192853505
$ws.Range("E4").Value = 0.000319350852
$ws.Range("I4").Formula = "=SUM(C4:H4)"

# J3/J4 share the same relative formula pattern -> fill as one range so the
# engine records them as a shared-formula group like the original edit did.
$ws.Range("J3:J4").Formula = "=I3*100*100"

# --- Row 5: mortality ratio ---
$ws.Range("I5").Value = "mortality"
$ws.Range("I5").Font.Bold = $true

$ws.Range("J5").Formula = "=(J4/J3)"
$ws.Range("J5").Style = "Percent"
$ws.Range("J5").NumberFormat = "0.00%"
$ws.Range("J5").Font.Size = 12
$ws.Range("J5").Font.ThemeColor = 1
$ws.Range("J5").Font.Bold = $true

# --- column J width (bestFit) ---
$ws.Columns.Item(10).ColumnWidth = 12.33203125

# --- sheet selection ---
$ws.Range("G12").Select()
